$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item(1)

# Add IDSheet after "data"
$idSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $dataSheet)
$idSheet.Name = "IDSheet"
$idSheet.Range("A1").Value = "KPI ID"
$idSheet.Range("A2").Value = "KC ID"
$idSheet.Range("B1").Value = "aZCyzqYa2aqEtf2945cna6"
$idSheet.Range("B2").Value = "524fc08b8a0e4d8d857dded88d5fb882"

# Ensure numbers in data sheet column A (rows 3-6) are stored as text
$dataSheet.Range("A3").NumberFormat = "@"
$dataSheet.Range("A3").Value = "'55"
$dataSheet.Range("A4").NumberFormat = "@"
$dataSheet.Range("A4").Value = "'22"
$dataSheet.Range("A5").NumberFormat = "@"
$dataSheet.Range("A5").Value = "'56"
$dataSheet.Range("A6").NumberFormat = "@"
$dataSheet.Range("A6").Value = "'24"

# Selections: IDSheet's selection sits on M34; data sheet stays the active tab,
# selection moves to A7
$idSheet.Range("M34").Select() | Out-Null
$dataSheet.Activate() | Out-Null
$dataSheet.Range("A7").Select() | Out-Null
